# "Add figures for 'Rectangles Everywhere!' puzzle."
#
# Slide 3 contains one top-level Group ("Group 6") holding the figure's
# rectangles, braces and dot-textboxes. We:
#   1. Ungroup it so the member shapes become top-level and editable.
#   2. Reposition/resize/resize-font the ". . ." textbox that sits between
#      the top-right rectangles.
#   3. Delete the three dot-textboxes that marked the "vertical" ellipsis
#      runs (they're replaced by rotated ". . ." textboxes below).
#   4. Add three new ". . ." textboxes (two rotated 90 degrees) to mark the
#      omitted rows/column of rectangles.
#   5. Re-select everything and re-group -- this is what actually produces
#      PowerPoint's renumbered "Group 3" (id 4), matching real PowerPoint
#      behavior when a group is ungrouped and regrouped after edits.
#
# NOTE: EMU->point conversion on this host truncates instead of rounding,
# so a tiny positive epsilon (well under 1 EMU) is added to each point
# value to land exactly on the target EMU after the point->EMU round trip.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$EMU = 12700.0
$EPS = 0.00002

function EmuPt([double]$emu) {
    return ($emu / $EMU) + $EPS
}

$g = $s.Shapes.Item(1)
$g.Ungroup() | Out-Null

# --- Locate member shapes by name (z-order may shift as we delete/add) ---
function FindShapeByName([string]$name) {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        if ($s.Shapes.Item($i).Name -eq $name) {
            return $s.Shapes.Item($i)
        }
    }
    return $null
}

# 2. Reposition / resize / font-size the ". . ." textbox ("TextBox 2").
$tb2 = FindShapeByName "TextBox 2"
$tb2.Left = EmuPt 5517928
$tb2.Top = EmuPt 1528897
$tb2.Width = EmuPt 1324303
$tb2.Height = EmuPt 400110
$tb2.TextFrame.TextRange.Font.Size = 20

# 3. Delete the three obsolete dot textboxes.
(FindShapeByName "TextBox 33").Delete() | Out-Null
(FindShapeByName "TextBox 34").Delete() | Out-Null
(FindShapeByName "TextBox 35").Delete() | Out-Null

# 4. Add the three new ". . ." textboxes.
$new1 = $s.Shapes.AddTextbox(1, (EmuPt 3579130), (EmuPt 2422275), (EmuPt 1324303), (EmuPt 400110))
$new1.Fill.Visible = 0
$new1.TextFrame.WordWrap = -1
$new1.TextFrame.AutoSize = 1
$new1.TextFrame.VerticalAnchor = 3
$new1.TextFrame.TextRange.Text = ".   .   ."
$new1.TextFrame.TextRange.Font.Size = 20
$new1.TextFrame.TextRange.Font.Bold = -1
$new1.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$new1.Rotation = 90

$new2 = $s.Shapes.AddTextbox(1, (EmuPt 7704440), (EmuPt 2419439), (EmuPt 1324303), (EmuPt 400110))
$new2.Fill.Visible = 0
$new2.TextFrame.WordWrap = -1
$new2.TextFrame.AutoSize = 1
$new2.TextFrame.VerticalAnchor = 3
$new2.TextFrame.TextRange.Text = ".   .   ."
$new2.TextFrame.TextRange.Font.Size = 20
$new2.TextFrame.TextRange.Font.Bold = -1
$new2.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$new2.Rotation = 90

$new3 = $s.Shapes.AddTextbox(1, (EmuPt 5517928), (EmuPt 3315657), (EmuPt 1324303), (EmuPt 400110))
$new3.Fill.Visible = 0
$new3.TextFrame.WordWrap = -1
$new3.TextFrame.AutoSize = 1
$new3.TextFrame.VerticalAnchor = 3
$new3.TextFrame.TextRange.Text = ".   .   ."
$new3.TextFrame.TextRange.Font.Size = 20
$new3.TextFrame.TextRange.Font.Bold = -1
$new3.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# 5. Re-select all member shapes (in z-order) and re-group them.
$all = $s.Shapes.Range((1..$s.Shapes.Count))
$all.Group() | Out-Null
